$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1802600.2
$ws.Range("I18").Value = 2250750.2
$ws.Range("K18").Value = 2250750.2
$ws.Range("M18").Value = -2250466.2

$ws.Range("H40").Value = 2731.375
$ws.Range("I40").Value = 2785.7144
$ws.Range("K40").Value = 2785.7144
$ws.Range("M40").Value = -2610.7144

$ws.Range("H95").Value = 40812
$ws.Range("J95").Value = 40812
$ws.Range("L95").Value = 40812
$ws.Range("N95").Value = -46304

$ws.Range("H100").Value = 6297.2
$ws.Range("I100").Value = 7622.5
$ws.Range("J100").Value = 996
$ws.Range("K100").Value = 7622.5
$ws.Range("L100").Value = 996
$ws.Range("M100").Value = -7081.5
$ws.Range("N100").Value = -2078

$ws.Range("H106").Value = 74502.57000000001
$ws.Range("I106").Value = 3337.3333
$ws.Range("K106").Value = 3337.3333
$ws.Range("M106").Value = -2706.3333

$ws.Range("H107").Value = 381.53333
$ws.Range("I107").Value = 434.33334
$ws.Range("J107").Value = 302.33334
$ws.Range("K107").Value = 434.33334
$ws.Range("L107").Value = 302.33334
$ws.Range("M107").Value = 1485.66666
$ws.Range("N107").Value = -4142.33334

$ws.Range("H116").Value = 17716236
$ws.Range("I116").Value = 10468285
$ws.Range("J116").Value = 37044104
$ws.Range("K116").Value = 10468285
$ws.Range("L116").Value = 37044104
$ws.Range("M116").Value = -10464843
$ws.Range("N116").Value = -37050988

$ws.Range("H135").Value = 797.75
$ws.Range("I135").Value = 797.75
$ws.Range("K135").Value = 7179.75
$ws.Range("M135").Value = -4644.75

$ws.Range("H137").Value = 3658.5334
$ws.Range("I137").Value = 2821.0908
$ws.Range("K137").Value = 8463.2724
$ws.Range("M137").Value = -5913.2724

$ws.Range("H138").Value = 2098.5386
$ws.Range("I138").Value = 1889.7037
$ws.Range("J138").Value = 2324.08
$ws.Range("K138").Value = 5669.1111
$ws.Range("L138").Value = 6972.24
$ws.Range("M138").Value = -529.1111000000001
$ws.Range("N138").Value = -17252.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8721.195
$ws.Range("I32").Value = 5285.0986
$ws.Range("K32").Value = 5285.0986
$ws.Range("M32").Value = -4998.0986

$ws.Range("H122").Value = 4095.0645
$ws.Range("I122").Value = 3698.4167
$ws.Range("J122").Value = 5455
$ws.Range("K122").Value = 11095.2501
$ws.Range("L122").Value = 16365
$ws.Range("M122").Value = -8645.250100000001
$ws.Range("N122").Value = -21265

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2447.8667
$ws.Range("I20").Value = 2746.3
$ws.Range("J20").Value = 1851
$ws.Range("K20").Value = 2746.3
$ws.Range("L20").Value = 1851
$ws.Range("M20").Value = -2499.3
$ws.Range("N20").Value = -2345

$ws.Range("H94").Value = 13891420
$ws.Range("I94").Value = 13891420
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 13891420
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -13890969
$ws.Range("N94").ClearContents()

$ws.Range("H134").Value = 4297.857
$ws.Range("I134").Value = 3653.6667
$ws.Range("J134").Value = 5457.4
$ws.Range("K134").Value = 10961.0001
$ws.Range("L134").Value = 16372.2
$ws.Range("M134").Value = -8426.000100000001
$ws.Range("N134").Value = -21442.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 250
$ws.Range("I3").Value = 250
$ws.Range("K3").Value = 250
$ws.Range("M3").Value = -137

$ws.Range("H10").Value = 107.75
$ws.Range("I10").Value = 135
$ws.Range("J10").Value = 26
$ws.Range("K10").Value = 135
$ws.Range("L10").Value = 26
$ws.Range("M10").Value = 4
$ws.Range("N10").Value = -304

$ws.Range("H16").Value = 3826.75
$ws.Range("I16").Value = 3083.818
$ws.Range("K16").Value = 3083.818
$ws.Range("M16").Value = -2796.818

$ws.Range("H39").Value = 2499.75
$ws.Range("I39").Value = 2499.75
$ws.Range("K39").Value = 2499.75
$ws.Range("M39").Value = -2108.75

$ws.Range("H44").Value = 9000
$ws.Range("I44").Value = 9000
$ws.Range("K44").Value = 9000
$ws.Range("M44").Value = -8558

$ws.Range("H49").Value = 2499.75
$ws.Range("I49").Value = 2499.75
$ws.Range("K49").Value = 2499.75
$ws.Range("M49").Value = -2317.75

$ws.Range("H105").Value = 780.6667
$ws.Range("I105").Value = 816.2727
$ws.Range("K105").Value = 816.2727
$ws.Range("M105").Value = 930.7273

$ws.Range("H113").Value = 3826.75
$ws.Range("I113").Value = 3083.818
$ws.Range("K113").Value = 3083.818
$ws.Range("M113").Value = -913.8180000000002

$ws.Range("H141").Value = 214178
$ws.Range("J141").Value = 226984
$ws.Range("L141").Value = 226984
$ws.Range("N141").Value = -237344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1172893
$ws.Range("I11").Value = 38676.074
$ws.Range("K11").Value = 116028.222
$ws.Range("M11").Value = -115888.222

$ws.Range("H38").Value = 251.85715
$ws.Range("I38").Value = 272.2
$ws.Range("J38").Value = 201
$ws.Range("K38").Value = 816.5999999999999
$ws.Range("L38").Value = 603
$ws.Range("M38").Value = -469.5999999999999

$ws.Range("I122").Value = 990
$ws.Range("J122").Value = 66667330
$ws.Range("K122").Value = 8910
$ws.Range("L122").Value = 600005970
$ws.Range("M122").Value = -6460
$ws.Range("N122").Value = -600010870

$ws.Range("H131").Value = 10026.074
$ws.Range("I131").Value = 1909.5
$ws.Range("J131").Value = 14800.529
$ws.Range("K131").Value = 5728.5
$ws.Range("L131").Value = 44401.587
$ws.Range("M131").Value = -688.5
$ws.Range("N131").Value = -54481.587

$ws.Range("H139").Value = 2787963
$ws.Range("J139").Value = 2812.5
$ws.Range("L139").Value = 8437.5
$ws.Range("N139").Value = -18717.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 7617.25
$ws.Range("J101").Value = 7617.25
$ws.Range("L101").Value = 7617.25
$ws.Range("N101").Value = -14107.25

$ws.Range("H102").Value = 2142.389
$ws.Range("I102").Value = 1249.5
$ws.Range("K102").Value = 1249.5
$ws.Range("M102").Value = 372.5

$ws.Range("H107").Value = 657.9
$ws.Range("I107").Value = 631
$ws.Range("K107").Value = 631
$ws.Range("M107").Value = 1289

$ws.Range("H113").Value = 4434.2
$ws.Range("I113").Value = 2227
$ws.Range("K113").Value = 2227
$ws.Range("M113").Value = -57

$ws.Range("H123").Value = 22074
$ws.Range("J123").Value = 22074
$ws.Range("L123").Value = 22074
$ws.Range("N123").Value = -26974

$ws.Range("H126").Value = 11821.533
$ws.Range("I126").Value = 22420.666
$ws.Range("K126").Value = 67261.99800000001
$ws.Range("M126").Value = -64791.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 361.8
$ws.Range("I22").Value = 358.75
$ws.Range("J22").Value = 374
$ws.Range("K22").Value = 358.75
$ws.Range("L22").Value = 374
$ws.Range("M22").Value = -63.75
$ws.Range("N22").Value = -964

$ws.Range("H27").Value = 361.8
$ws.Range("I27").Value = 358.75
$ws.Range("J27").Value = 374
$ws.Range("K27").Value = 358.75
$ws.Range("L27").Value = 374
$ws.Range("M27").Value = -251.75
$ws.Range("N27").Value = -588

$ws.Range("H46").Value = 2562.7144
$ws.Range("I46").Value = 2033
$ws.Range("K46").Value = 2033
$ws.Range("M46").Value = -1845

$ws.Range("H93").Value = 83335490
$ws.Range("I93").Value = 83335490
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 83335490
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -83334242
$ws.Range("N93").ClearContents()

$ws.Range("H100").Value = 4458.077
$ws.Range("I100").Value = 4217.222
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 4217.222
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -3676.222
$ws.Range("N100").Value = -6082

$ws.Range("H136").Value = 4466.724
$ws.Range("I136").Value = 4295.8037
$ws.Range("K136").Value = 12887.4111
$ws.Range("M136").Value = -10337.4111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 10001181
$ws.Range("J100").Value = 965.6667
$ws.Range("L100").Value = 1931.3334
$ws.Range("N100").Value = -3013.3334

$ws.Range("H122").Value = 2788
$ws.Range("I122").Value = 3333.3333
$ws.Range("J122").Value = 2490.5454
$ws.Range("K122").Value = 9999.999899999999
$ws.Range("L122").Value = 7471.6362
$ws.Range("M122").Value = -7549.999899999999
$ws.Range("N122").Value = -12371.6362

$ws.Range("H132").Value = 28996.625
$ws.Range("I132").Value = 11742.363
$ws.Range("K132").Value = 35227.089
$ws.Range("M132").Value = -32697.089

$ws.Range("H140").Value = 76985.28999999999
$ws.Range("J140").Value = 76985.28999999999
$ws.Range("L140").Value = 76985.28999999999
$ws.Range("N140").Value = -87345.28999999999

$ws.Range("H141").Value = 89999
$ws.Range("J141").Value = 89999
$ws.Range("L141").Value = 89999
$ws.Range("N141").Value = -100359
